$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H2").Value = [double]"0.1092515172334684"
$ws.Range("B3").Value = [double]"0.0003786836924110576"
$ws.Range("C3").Value = [double]"0.0005234618411666489"
$ws.Range("D3").Value = [double]"3.051430761016436"
$ws.Range("E3").Value = [double]"0.0229775216029635"
$ws.Range("F3").Value = [double]"-0.0006472865170047952"
$ws.Range("G3").Value = [double]"0.00140465390182691"
$ws.Range("H3").Value = [double]"0.1096302009258794"
$ws.Range("B4").Value = [double]"0.003005056176465428"
$ws.Range("C4").Value = [double]"0.000716442962551953"
$ws.Range("D4").Value = [double]"6.358182644091029"
$ws.Range("E4").Value = [double]"0.01609570926523277"
$ws.Range("F4").Value = [double]"0.001600848385325376"
$ws.Range("G4").Value = [double]"0.00440926396760548"
$ws.Range("H4").Value = [double]"0.1122565734099338"
$ws.Range("B5").Value = [double]"0.01573271205740587"
$ws.Range("C5").Value = [double]"0.001374835634167393"
$ws.Range("D5").Value = [double]"11.34247894681815"
$ws.Range("E5").Value = [double]"0.08581123217433899"
$ws.Range("F5").Value = [double]"0.01303807334480385"
$ws.Range("G5").Value = [double]"0.01842735077000788"
$ws.Range("H5").Value = [double]"0.1249842292908742"
$ws.Range("B6").Value = [double]"0.02437826814118802"
$ws.Range("C6").Value = [double]"0.004721501975969864"
$ws.Range("D6").Value = [double]"9.644055900210581"
$ws.Range("E6").Value = [double]"0.2153419071376526"
$ws.Range("F6").Value = [double]"0.01512425957923483"
$ws.Range("G6").Value = [double]"0.03363227670314121"
$ws.Range("H6").Value = [double]"0.1336297853746564"
$ws.Range("B7").Value = [double]"0.02616946103183122"
$ws.Range("C7").Value = [double]"0.004263896006435288"
$ws.Range("D7").Value = [double]"8.112483971818488"
$ws.Range("E7").Value = [double]"0.07125234345880059"
$ws.Range("F7").Value = [double]"0.01781235084809405"
$ws.Range("G7").Value = [double]"0.03452657121556838"
$ws.Range("H7").Value = [double]"0.1354209782652996"
$ws.Range("B8").Value = [double]"0.02939336363868053"
$ws.Range("C8").Value = [double]"0.004065125816084954"
$ws.Range("D8").Value = [double]"9.955970666601123"
$ws.Range("E8").Value = [double]"0.04948282583821543"
$ws.Range("F8").Value = [double]"0.02142583922326095"
$ws.Range("G8").Value = [double]"0.03736088805410012"
$ws.Range("H8").Value = [double]"0.1386448808721489"
$ws.Range("B9").Value = [double]"0.03654309030418693"
$ws.Range("C9").Value = [double]"0.002641373709231868"
$ws.Range("D9").Value = [double]"10.08081428379742"
$ws.Range("E9").Value = [double]"0.008602939053818722"
$ws.Range("F9").Value = [double]"0.03136607763740564"
$ws.Range("G9").Value = [double]"0.04172010297096822"
$ws.Range("H9").Value = [double]"0.1457946075376553"
$ws.Range("B10").Value = [double]"-0.1092515172334684"
$ws.Range("C10").Value = [double]"0.0004417530203152026"
$ws.Range("D10").Value = [double]"-258.5107752082333"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.1101173406010237"
$ws.Range("G10").Value = [double]"-0.108385693865913"
$ws.Range("B11").Value = [double]"-0.04923659764432254"
$ws.Range("C11").Value = [double]"0.0005032604742101206"
$ws.Range("D11").Value = [double]"-101.8510236846971"
$ws.Range("E11").Value = [double]"0"
$ws.Range("F11").Value = [double]"-0.05022297394184667"
$ws.Range("G11").Value = [double]"-0.04825022134679842"
$ws.Range("H11").Value = [double]"0.06001491958914582"
$ws.Range("B12").Value = [double]"-0.04358314019471648"
$ws.Range("C12").Value = [double]"0.0004922985833009629"
$ws.Range("D12").Value = [double]"-92.8253905059364"
$ws.Range("E12").Value = [double]"1.003617431177309e-300"
$ws.Range("F12").Value = [double]"-0.04454803153826932"
$ws.Range("G12").Value = [double]"-0.04261824885116364"
$ws.Range("H12").Value = [double]"0.06566837703875189"
$ws.Range("B13").Value = [double]"-0.03525661579933682"
$ws.Range("C13").Value = [double]"0.0004836667498322539"
$ws.Range("D13").Value = [double]"-75.33782534724892"
$ws.Range("E13").Value = [double]"1.585710338892861e-308"
$ws.Range("F13").Value = [double]"-0.03620458889803842"
$ws.Range("G13").Value = [double]"-0.03430864270063522"
$ws.Range("H13").Value = [double]"0.07399490143413154"
$ws.Range("B14").Value = [double]"-0.02917869067366006"
$ws.Range("C14").Value = [double]"0.000472303102981849"
$ws.Range("D14").Value = [double]"-62.94586745426815"
$ws.Range("E14").Value = [double]"7.090887278709099e-201"
$ws.Range("F14").Value = [double]"-0.03010439128024484"
$ws.Range("G14").Value = [double]"-0.02825299006707527"
$ws.Range("H14").Value = [double]"0.08007282655980831"
$ws.Range("B15").Value = [double]"-0.02430368739204947"
$ws.Range("C15").Value = [double]"0.0004533927914866108"
$ws.Range("D15").Value = [double]"-54.62890543909047"
$ws.Range("E15").Value = [double]"6.182682758706705e-86"
$ws.Range("F15").Value = [double]"-0.02519232432669894"
$ws.Range("G15").Value = [double]"-0.0234150504574"
$ws.Range("H15").Value = [double]"0.0849478298414189"
$ws.Range("B16").Value = [double]"-0.02198077672140187"
$ws.Range("C16").Value = [double]"0.0004452368967848326"
$ws.Range("D16").Value = [double]"-50.43381914662623"
$ws.Range("E16").Value = [double]"1.791947137627677e-73"
$ws.Range("F16").Value = [double]"-0.02285342842974601"
$ws.Range("G16").Value = [double]"-0.02110812501305772"
$ws.Range("H16").Value = [double]"0.0872707405120665"
$ws.Range("B17").Value = [double]"-0.02041562309456206"
$ws.Range("C17").Value = [double]"0.0004589752218225486"
$ws.Range("D17").Value = [double]"-45.68227010951593"
$ws.Range("E17").Value = [double]"1.339812432967346e-33"
$ws.Range("F17").Value = [double]"-0.02131520148114452"
$ws.Range("G17").Value = [double]"-0.01951604470797961"
$ws.Range("H17").Value = [double]"0.0888358941389063"
$ws.Range("B18").Value = [double]"-0.01828536644378089"
$ws.Range("C18").Value = [double]"0.0004573163500083862"
$ws.Range("D18").Value = [double]"-42.14415701262025"
$ws.Range("E18").Value = [double]"9.432821093385412e-17"
$ws.Range("F18").Value = [double]"-0.0191816934424996"
$ws.Range("G18").Value = [double]"-0.01738903944506218"
$ws.Range("H18").Value = [double]"0.09096615078968748"
$ws.Range("B19").Value = [double]"-0.01451477120554883"
$ws.Range("C19").Value = [double]"0.0004466295061027155"
$ws.Range("D19").Value = [double]"-33.77353424832028"
$ws.Range("E19").Value = [double]"0.0258078019153806"
$ws.Range("F19").Value = [double]"-0.01539015226009925"
$ws.Range("G19").Value = [double]"-0.01363939015099841"
$ws.Range("H19").Value = [double]"0.09473674602791952"
$ws.Range("B20").Value = [double]"-0.01099169640886734"
$ws.Range("C20").Value = [double]"0.0004579375002018215"
$ws.Range("D20").Value = [double]"-25.37749700439542"
$ws.Range("E20").Value = [double]"0.009463938449336297"
$ws.Range("F20").Value = [double]"-0.01188924084957786"
$ws.Range("G20").Value = [double]"-0.01009415196815682"
$ws.Range("H20").Value = [double]"0.09825982082460102"
$ws.Range("B21").Value = [double]"-0.00838126266252689"
$ws.Range("C21").Value = [double]"0.0004547068561464612"
$ws.Range("D21").Value = [double]"-19.66723917847885"
$ws.Range("E21").Value = [double]"3.800522286159876e-07"
$ws.Range("F21").Value = [double]"-0.00927247509842789"
$ws.Range("G21").Value = [double]"-0.007490050226625889"
$ws.Range("H21").Value = [double]"0.1008702545709415"
$ws.Range("B22").Value = [double]"-0.006516578438563062"
$ws.Range("C22").Value = [double]"0.0004500718100747738"
$ws.Range("D22").Value = [double]"-16.06690342537115"
$ws.Range("E22").Value = [double]"2.072965561200962e-06"
$ws.Range("F22").Value = [double]"-0.007398706292282006"
$ws.Range("G22").Value = [double]"-0.005634450584844118"
$ws.Range("H22").Value = [double]"0.1027349387949053"
$ws.Range("B23").Value = [double]"-0.004834377299731717"
$ws.Range("C23").Value = [double]"0.0004457097487561667"
$ws.Range("D23").Value = [double]"-12.45685317456049"
$ws.Range("E23").Value = [double]"0.06499528193420305"
$ws.Range("F23").Value = [double]"-0.005707955637653282"
$ws.Range("G23").Value = [double]"-0.00396079896181015"
$ws.Range("H23").Value = [double]"0.1044171399337366"
$ws.Range("B24").Value = [double]"-0.002853841335796048"
$ws.Range("C24").Value = [double]"0.0004317966924833947"
$ws.Range("D24").Value = [double]"-7.066509358383073"
$ws.Range("E24").Value = [double]"0.04722473584311705"
$ws.Range("F24").Value = [double]"-0.00370015047167107"
$ws.Range("G24").Value = [double]"-0.002007532199921025"
$ws.Range("H24").Value = [double]"0.1063976758976723"
$ws.Range("B25").Value = [double]"-0.00124576507936059"
$ws.Range("C25").Value = [double]"0.000410807427333105"
$ws.Range("D25").Value = [double]"-3.349684285051641"
$ws.Range("E25").Value = [double]"0.07982811500891104"
$ws.Range("F25").Value = [double]"-0.002050935879165034"
$ws.Range("G25").Value = [double]"-0.000440594279556145"
$ws.Range("H25").Value = [double]"0.1080057521541078"
$ws.Range("B26").Value = [double]"0.07860134012120554"
$ws.Range("C26").Value = [double]"0.002720485738972481"
$ws.Range("D26").Value = [double]"54.69272957287313"
$ws.Range("E26").Value = [double]"1.530880959535897e-07"
$ws.Range("F26").Value = [double]"0.07326925762690703"
$ws.Range("G26").Value = [double]"0.08393342261550404"
$ws.Range("H26").Value = [double]"0.1878528573546739"
